$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update week 2 (row 8) hours from "~1" to "~2.5"
$ws.Range("E8").Value = "~2.5"

# Fill in week 3 (row 9) hours and activities
$ws.Range("E9").Value = 2.5
$ws.Range("F9").Value = "Backlog template updates, planning, time sheet"

# Update the active selection to match the authored state
$ws.Range("E8").Select()
